$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: update values, clear C2
$ws.Range("B2").Value = 93.391657235873353
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 94.440979878641414
$ws.Range("E2").Value = 94.450505279159486

# Row 3: clear B3, update C3, add D3, update E3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 97.000469506419734
$ws.Range("D3").Value = 95.128685951079902
$ws.Range("E3").Value = 94.700599725404615

# Selection / view tweak
$ws.Range("B1:E3").Select()
